$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column A to fit the new, longer label that gets added below ---
$ws.Columns.Item(1).ColumnWidth = 22.28

# --- Row 2: keep A2 ("Cumplimiento de pago"); redistribute its score cells ---
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 1
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 5

# --- Row 3 (new): "Cuidado del Inmueble" ---
$ws.Range("A3").Value = "Cuidado del Inmueble"
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = 5

# --- Row 4 (new): "Responsabilidad General" ---
$ws.Range("A4").Value = "Responsabilidad General"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 4

# --- Extend each chart series to cover the new rows (A2:A4 / x2:x4) ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES('Sheet1'!B1,Sheet1!A2:A4,Sheet1!B2:B4,1)"

$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES('Sheet1'!C1,Sheet1!A2:A4,Sheet1!C2:C4,2)"

$s3 = $chart.SeriesCollection().Item(3)
$s3.Formula = "=SERIES('Sheet1'!D1,Sheet1!A2:A4,Sheet1!D2:D4,3)"

$s4 = $chart.SeriesCollection().Item(4)
$s4.Formula = "=SERIES('Sheet1'!E1,Sheet1!A2:A4,Sheet1!E2:E4,4)"

$s5 = $chart.SeriesCollection().Item(5)
$s5.Formula = "=SERIES('Sheet1'!F1,Sheet1!A2:A4,Sheet1!F2:F4,5)"

# --- Shrink the chart's right edge a bit (its "to" column offset) ---
# Read the live width (reflects the column-A resize above) and shift it by
# the exact EMU delta so the anchor's colOff moves from 295275 to 171450.
$currentWidth = $co.Width
$newWidth = $currentWidth - (295275 - 171450) / 12700.0
$co.Width = $newWidth
